$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Month" column (B) with header + two month names, matching the
# existing A1:A3 header/value layout.
$ws.Range("B1").Value = "Month"
$ws.Range("B2").Value = "Helmikuu"
$ws.Range("B3").Value = "Maaliskuu"

# Column A was manually narrowed (no longer "best fit"); column B sized
# to comfortably fit the month names.
$ws.Range("A1").ColumnWidth = 43.8333333333333
$ws.Range("B1").ColumnWidth = 9

# Selection moves to the newly entered cell.
[void]$ws.Range("B3").Select()
